# u matrix: electricity use from gas boiler put to zero
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Row 4 = "Electricity" activity; set electricity consumption of
# gas-fired appliances (Gas boiler for Heating, Gas boiler for Hot
# Sanitary Water, Gas Stove for Cooking) to zero.
$ws.Range("H4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("P4").Value = 0
